$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM (ligand-receptor) values for the Proc-Tek sheet
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.347467
$ws.Range("H2").Value = 1.042401
$ws.Range("I2").Value = 0.864291305025554
$ws.Range("J2").Value = 0.864291305025554
$ws.Range("M2").Value = 71.05094633333333
$ws.Range("N2").Value = 213.152839
$ws.Range("O2").Value = 0.8240565632932695
$ws.Range("P2").Value = 0.8240565632932696
$ws.Range("Q2").Value = 24.68785916960433
$ws.Range("R2").Value = 222.190732526439
$ws.Range("S2").Value = 0.7122249225036129
$ws.Range("T2").Value = 0.712224922503613
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.347467
$ws.Range("H3").Value = 1.042401
$ws.Range("I3").Value = 0.864291305025554
$ws.Range("J3").Value = 0.864291305025554
$ws.Range("O3").Value = 0.1323102827659759
$ws.Range("P3").Value = 0.132310282765976
$ws.Range("Q3").Value = 3.963875506995333
$ws.Range("R3").Value = 35.674879562958
$ws.Range("S3").Value = 0.1143546269601054
$ws.Range("T3").Value = 0.1143546269601054
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.347467
$ws.Range("H4").Value = 1.042401
$ws.Range("I4").Value = 0.864291305025554
$ws.Range("J4").Value = 0.864291305025554
$ws.Range("M4").Value = 3.762092333333333
$ws.Range("N4").Value = 11.286277
$ws.Range("O4").Value = 0.04363315394075456
$ws.Range("P4").Value = 0.04363315394075455
$ws.Range("Q4").Value = 1.307202936786333
$ws.Range("R4").Value = 11.764826431077
$ws.Range("S4").Value = 0.03771175556183565
$ws.Range("T4").Value = 0.03771175556183565
$ws.Range("I5").Value = 0.135708694974446
$ws.Range("J5").Value = 0.135708694974446
$ws.Range("M5").Value = 71.05094633333333
$ws.Range("N5").Value = 213.152839
$ws.Range("O5").Value = 0.8240565632932695
$ws.Range("P5").Value = 0.8240565632932696
$ws.Range("Q5").Value = 3.876421213702777
$ws.Range("R5").Value = 34.887790923325
$ws.Range("S5").Value = 0.1118316407896566
$ws.Range("T5").Value = 0.1118316407896566
$ws.Range("I6").Value = 0.135708694974446
$ws.Range("J6").Value = 0.135708694974446
$ws.Range("O6").Value = 0.1323102827659759
$ws.Range("P6").Value = 0.132310282765976
$ws.Range("S6").Value = 0.01795565580587053
$ws.Range("T6").Value = 0.01795565580587053
$ws.Range("I7").Value = 0.135708694974446
$ws.Range("J7").Value = 0.135708694974446
$ws.Range("O7").Value = 0.04363315394075456
$ws.Range("P7").Value = 0.04363315394075455
$ws.Range("Q7").Value = 0.2052534875527777
$ws.Range("R7").Value = 1.847281387975
$ws.Range("S7").Value = 0.005921398378918909
$ws.Range("T7").Value = 0.005921398378918908